# Cambios en datos entidades para que funcione con el nuevo formato
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "concepto" (D) and "razonEntidad" (E) values for ICBF (row 20) and ESAP (row 22)
# so columns line up with the new expected format.
$d20 = $ws.Range("D20").Value2
$e20 = $ws.Range("E20").Value2
$ws.Range("D20").Value2 = $e20
$ws.Range("E20").Value2 = $d20

$d22 = $ws.Range("D22").Value2
$e22 = $ws.Range("E22").Value2
$ws.Range("D22").Value2 = $e22
$ws.Range("E22").Value2 = $d22

# Update the active selection left in the sheet view.
$ws.Range("E21").Select()
